$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.077.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.51%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.268.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.12%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.53%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.641"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.60%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.38"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.90%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.12%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.659"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +15.83%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.23%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0976"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.94%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "59.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.91%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.74%  "

# Row 14
$ws.Range("E14").Value = "  +0.98%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.606.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.24%  "

# Row 16
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "Chainlink"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.91"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.02%  "

# Row 17
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "Polygon"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.890"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.04%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.255.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.94%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.953.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.54%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0984"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.33%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.57%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.75%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.64%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.95%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.02%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.71%  "

# Row 27
$ws.Range("E27").Value = "  +0.01%  "

# Row 28
$ws.Range("E28").Value = "  -1.99%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.15%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.32%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "168.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.27%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.37%  "

# Row 33
$ws.Range("E33").Value = "  +11.81%  "

# Row 34
$ws.Range("E34").Value = "  +13.80%  "

# Row 35
$ws.Range("E35").Value = "  +5.52%  "

# Row 36
$ws.Range("E36").Value = "  +2.90%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "28.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.56%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.50%  "

# Row 39
$ws.Range("E39").Value = "  +1.31%  "

# Row 40
$ws.Range("E40").Value = "  +8.74%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.71%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.05%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.85"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.92%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.03%  "

# Row 45
$ws.Range("E45").Value = "  +0.27%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.203"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.74%  "

# Row 47
$ws.Range("E47").Value = "  +5.66%  "

# Row 48
$ws.Range("E48").Value = "  +2.26%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.98%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.24%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.43%  "

